$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells touched by this update (Price / Volume(1h) columns for several coins).
# The sheet stores these as plain text (e.g. "277.48", "0.94%") rather than
# numbers/percentages, so force text format on each cell first - otherwise
# Excel would "smart" parse a string like "0.94%" into a percentage number.
$changedCells = @("D2", "E2", "D3", "E3", "D4", "E4", "D5", "E5", "D6", "E6", "E7", "D8", "E8", "D9", "E9", "D10", "E10", "D11", "E11", "D12", "E12", "D13", "E13", "D14", "E14", "D15", "E15", "D16", "E16", "D17", "D18", "E18", "D19", "E19", "D21", "E21", "D22", "E22", "D23", "E23", "D24", "E24", "E26", "D27", "E27", "E28", "E29", "E40", "D41", "E41", "E42", "D43", "E43", "D44", "E44", "D45", "E45", "D46", "E46", "D47", "E47")
foreach ($cellRef in $changedCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$ws.Range("D2").Value = "277.48"
$ws.Range("E2").Value = "0.94%"
$ws.Range("D3").Value = "27.32"
$ws.Range("E3").Value = "2.18%"
$ws.Range("D4").Value = "4.860"
$ws.Range("E4").Value = "-0.73%"
$ws.Range("D5").Value = "0.06418"
$ws.Range("E5").Value = "1.36%"
$ws.Range("D6").Value = "6.954"
$ws.Range("E6").Value = "1.31%"
$ws.Range("E7").Value = "-6.67%"
$ws.Range("D8").Value = "0.8753"
$ws.Range("E8").Value = "0.69%"
$ws.Range("D9").Value = "0.1542"
$ws.Range("E9").Value = "-10.41%"
$ws.Range("D10").Value = "0.05172"
$ws.Range("E10").Value = "2.80%"
$ws.Range("D11").Value = "0.07431"
$ws.Range("E11").Value = "0.17%"
$ws.Range("D12").Value = "0.02956"
$ws.Range("E12").Value = "-0.17%"
$ws.Range("D13").Value = "0.08983"
$ws.Range("E13").Value = "-0.43%"
$ws.Range("D14").Value = "0.001569"
$ws.Range("E14").Value = "-0.35%"
$ws.Range("D15").Value = "0.0006373"
$ws.Range("E15").Value = "0.55%"
$ws.Range("D16").Value = "0.006122"
$ws.Range("E16").Value = "5.88%"
$ws.Range("D17").Value = "3.479"
$ws.Range("D18").Value = "3.308"
$ws.Range("E18").Value = "-0.23%"
$ws.Range("D19").Value = "2.274"
$ws.Range("E19").Value = "-0.43%"
$ws.Range("D21").Value = "0.1348"
$ws.Range("E21").Value = "1.02%"
$ws.Range("D22").Value = "3.903"
$ws.Range("E22").Value = "-0.40%"
$ws.Range("D23").Value = "0.04418"
$ws.Range("E23").Value = "1.27%"
$ws.Range("D24").Value = "0.1500"
$ws.Range("E24").Value = "8.66%"
$ws.Range("E26").Value = "-0.06%"
$ws.Range("D27").Value = "0.003870"
$ws.Range("E27").Value = "-8.93%"
$ws.Range("E28").Value = "8.22%"
$ws.Range("E29").Value = "15.06%"
$ws.Range("E40").Value = "2.24%"
$ws.Range("D41").Value = "0.006766"
$ws.Range("E41").Value = "0.00%"
$ws.Range("E42").Value = "0.54%"
$ws.Range("D43").Value = "0.001969"
$ws.Range("E43").Value = "-8.90%"
$ws.Range("D44").Value = "0.01148"
$ws.Range("E44").Value = "6.75%"
$ws.Range("D45").Value = "0.00005305"
$ws.Range("E45").Value = "0.28%"
$ws.Range("D46").Value = "1.687"
$ws.Range("E46").Value = "13.27%"
$ws.Range("D47").Value = "0.01852"
$ws.Range("E47").Value = "-11.86%"

# Restore the default "Normal" style so no stray number-format style lingers
# on cells that did not carry one before this edit.
foreach ($cellRef in $changedCells) {
    $ws.Range($cellRef).Style = "Normal"
}

